$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Expanded dataset")
$ws3 = $wb.Worksheets.Item("Adam optimizer")

# ------------------------------------------------------------------
# Defined names (new named ranges for the Adam-optimizer batch-norm
# learning curves), scoped to the "Adam optimizer" sheet.
# ------------------------------------------------------------------
$ws3.Names.Add("learning_curve_adam_batch_norm", "='Adam optimizer'!`$J`$4:`$K`$14") | Out-Null
$ws3.Names.Add("learning_curve_adam_batch_norm_7000", "='Adam optimizer'!`$L`$4:`$M`$14") | Out-Null

# ------------------------------------------------------------------
# Sheet2 "Expanded dataset": updated title + two new columns (J, K)
# holding the "batch norm" learning curve imported from
# learning_curve_adam_batch_norm.csv.
# ------------------------------------------------------------------
$ws2.Range("B2").Value2 = "Comparison of hyperparameters for the feed-forward neural network (lr = 0.03, expanded dataset with rotated and flipped images, minibatches = 7000)"

$ws2.Range("J4").Value2 = "train batch norm weight_decay: 0.0, dropout: 0.0"
$ws2.Range("K4").Value2 = "test batch norm weight_decay: 0.0, dropout: 0.0"

$ws2.Range("J5").Value2 = 0.91931746031745998
$ws2.Range("K5").Value2 = 0.85133333333333305
$ws2.Range("J6").Value2 = 0.94926984126984104
$ws2.Range("K6").Value2 = 0.84688888888888803
$ws2.Range("J7").Value2 = 0.96961904761904705
$ws2.Range("K7").Value2 = 0.84588888888888802
$ws2.Range("J8").Value2 = 0.98088888888888803
$ws2.Range("K8").Value2 = 0.84599999999999997
$ws2.Range("J9").Value2 = 0.98630158730158701
$ws2.Range("K9").Value2 = 0.846444444444444
$ws2.Range("J10").Value2 = 0.99047619047618995
$ws2.Range("K10").Value2 = 0.84044444444444399
$ws2.Range("J11").Value2 = 0.99166666666666603
$ws2.Range("K11").Value2 = 0.83855555555555505
$ws2.Range("J12").Value2 = 0.99296825396825295
$ws2.Range("K12").Value2 = 0.83711111111111103
$ws2.Range("J13").Value2 = 0.99631746031746005
$ws2.Range("K13").Value2 = 0.84333333333333305
$ws2.Range("J14").Value2 = 0.99653968253968195
$ws2.Range("K14").Value2 = 0.83911111111111103

# New column widths on sheet2 (cols J=10, K=11)
$ws2.Columns.Item(10).ColumnWidth = 29.6
$ws2.Columns.Item(11).ColumnWidth = 21.3

# ------------------------------------------------------------------
# Sheet3 "Adam optimizer": four new columns (J,K,L,M) holding the
# "batch norm" learning curves (minibatches==1000 and ==7000 variants)
# imported from the corresponding csv files.
# ------------------------------------------------------------------
$ws3.Range("J4").Value2 = "train batch norm weight_decay: 0.0, dropout: 0.0"
$ws3.Range("K4").Value2 = "test batch norm weight_decay: 0.0, dropout: 0.0"
$ws3.Range("L4").Value2 = "train batch norm weight_decay: 0.0, dropout: 0.0 mini batches = 7000"
$ws3.Range("M4").Value2 = "test batch norm weight_decay: 0.0, dropout: 0.0, mini batches = 7000"

$ws3.Range("J5").Value2 = 0.87326984126984097
$ws3.Range("K5").Value2 = 0.83355555555555505
$ws3.Range("L5").Value2 = 0.91695238095238096
$ws3.Range("M5").Value2 = 0.84866666666666601

$ws3.Range("J6").Value2 = 0.927539682539682
$ws3.Range("K6").Value2 = 0.83911111111111103
$ws3.Range("L6").Value2 = 0.95199999999999996
$ws3.Range("M6").Value2 = 0.845444444444444

$ws3.Range("J7").Value2 = 0.96244444444444399
$ws3.Range("K7").Value2 = 0.83533333333333304
$ws3.Range("L7").Value2 = 0.97047619047619005
$ws3.Range("M7").Value2 = 0.84566666666666601

$ws3.Range("J8").Value2 = 0.97847619047619006
$ws3.Range("K8").Value2 = 0.83422222222222198
$ws3.Range("L8").Value2 = 0.98120634920634897
$ws3.Range("M8").Value2 = 0.846444444444444

$ws3.Range("J9").Value2 = 0.986507936507936
$ws3.Range("K9").Value2 = 0.82799999999999996
$ws3.Range("L9").Value2 = 0.98682539682539605
$ws3.Range("M9").Value2 = 0.84233333333333305

$ws3.Range("J10").Value2 = 0.98969841269841197
$ws3.Range("K10").Value2 = 0.82522222222222197
$ws3.Range("L10").Value2 = 0.990079365079365
$ws3.Range("M10").Value2 = 0.842444444444444

$ws3.Range("J11").Value2 = 0.99433333333333296
$ws3.Range("K11").Value2 = 0.82788888888888801
$ws3.Range("L11").Value2 = 0.99406349206349198
$ws3.Range("M11").Value2 = 0.84266666666666601

$ws3.Range("J12").Value2 = 0.99676190476190396
$ws3.Range("K12").Value2 = 0.82377777777777705
$ws3.Range("L12").Value2 = 0.99553968253968195
$ws3.Range("M12").Value2 = 0.84033333333333304

$ws3.Range("J13").Value2 = 0.99798412698412697
$ws3.Range("K13").Value2 = 0.82355555555555504
$ws3.Range("L13").Value2 = 0.99461904761904696
$ws3.Range("M13").Value2 = 0.837666666666666

$ws3.Range("J14").Value2 = 0.99860317460317405
$ws3.Range("K14").Value2 = 0.82033333333333303
$ws3.Range("L14").Value2 = 0.99669841269841197
$ws3.Range("M14").Value2 = 0.83633333333333304

# New column widths on sheet3 (cols J=10, K=11, L=12, M=13)
$ws3.Columns.Item(10).ColumnWidth = 33.17
$ws3.Columns.Item(11).ColumnWidth = 40.8
$ws3.Columns.Item(12).ColumnWidth = 16.17
$ws3.Columns.Item(13).ColumnWidth = 19.97

# ------------------------------------------------------------------
# Selections (the author's final cursor position on each sheet).
# ------------------------------------------------------------------
$ws2.Range("B3").Select() | Out-Null
$ws3.Activate() | Out-Null
$ws3.Range("L6").Select() | Out-Null
